# Generate Report for Archive
# - Update the "Status" value shown for the two tracked files from
#   "Ready for handoff" to "In Translation" on all three sheets.
# - Narrow the corresponding "Status" column(s) to match the new,
#   shorter text (re-autofit of the column width).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: column C (Status) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: column C (Status) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
